$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "phone" detail row
$ws.Range("A3").Value = "phone"
$ws.Range("B3").Value = 6097364896

# Mirror the bestFit width captured in the diff for column B
$ws.Columns.Item(2).ColumnWidth = 9.92

# Update the active selection to match the authored state
$ws.Range("B3").Select() | Out-Null
